# A new weekly price observation is inserted above the existing row 152,
# pushing the former rows 152-207 down to 153-208 (dimension grows to R208).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(152).Insert()

$ws.Range("A152").Value = 11
$ws.Range("B152").Value = "Vega Monumental Concepción"
$ws.Range("C152").Value = "Bíobío"
$ws.Range("D152").Value = 45009
$ws.Range("E152").Value = 8
$ws.Range("F152").Value = 100112032
$ws.Range("G152").Value = "Zapallo italiano"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 100
$ws.Range("K152").Value = 7000
$ws.Range("L152").Value = 7500
$ws.Range("M152").Value = 7250
$ws.Range("N152").Value = "`$/caja 50 unidades"
$ws.Range("O152").Value = "Región de O'Higgins"
$ws.Range("P152").Value = 145
$ws.Range("Q152").Value = 50
$ws.Range("R152").Value = "Hortaliza"
